# Daily attendance processing - 2026-01-12 17:38:25
# Normalize the "Recorded By" (column G) text so the user email is listed
# before "System" instead of after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
